$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Item ID 5503)
$ws.Range("H5").Value = 260.33334
$ws.Range("I5").Value = 190
$ws.Range("J5").Value = 401
$ws.Range("K5").Value = 190
$ws.Range("L5").Value = 401
$ws.Range("M5").Value = -75
$ws.Range("N5").Value = -631

# Row 28 (Item ID 27772)
$ws.Range("H28").Value = 766.5
$ws.Range("I28").Value = 787.875
$ws.Range("J28").Value = 723.75
$ws.Range("K28").Value = 787.875
$ws.Range("L28").Value = 723.75
$ws.Range("M28").Value = -302.875
$ws.Range("N28").Value = -1693.75

# Row 62 (Item ID 27781)
$ws.Range("H62").Value = 940.5
$ws.Range("I62").Value = 875
$ws.Range("J62").Value = 1006
$ws.Range("K62").Value = 875
$ws.Range("L62").Value = 1006
$ws.Range("M62").Value = -251
$ws.Range("N62").Value = -2254

# Row 65 (Item ID 27781)
$ws.Range("H65").Value = 940.5
$ws.Range("I65").Value = 875
$ws.Range("J65").Value = 1006
$ws.Range("K65").Value = 4375
$ws.Range("L65").Value = 5030
$ws.Range("M65").Value = -1255
$ws.Range("N65").Value = -11270

# Row 106 (Item ID 19903)
$ws.Range("H106").Value = 266672460
$ws.Range("I106").Value = 111117450
$ws.Range("J106").Value = 500005000
$ws.Range("K106").Value = 111117450
$ws.Range("L106").Value = 500005000
$ws.Range("M106").Value = -111116819
$ws.Range("N106").Value = -500006262

# Row 129 (Item ID 36115)
$ws.Range("H129").Value = 1040.0435
$ws.Range("J129").Value = 1066.127
$ws.Range("L129").Value = 3198.381
$ws.Range("N129").Value = -13198.381

# Row 132 (Item ID 44049)
$ws.Range("H132").Value = 2077
$ws.Range("I132").Value = 1307.7778
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 3923.3334
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -1393.3334
$ws.Range("N132").Value = -32060

# Row 137 (Item ID 44013)
$ws.Range("H137").Value = 1863.2858
$ws.Range("I137").Value = 1239.4166
$ws.Range("J137").Value = 2695.111
$ws.Range("K137").Value = 3718.2498
$ws.Range("L137").Value = 8085.333
$ws.Range("M137").Value = -1168.2498
$ws.Range("N137").Value = -13185.333

$ws = $wb.Worksheets.Item("ARM")
# Row 110 (Item ID 27708)
$ws.Range("H110").Value = 751.375
$ws.Range("I110").Value = 668.5
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 668.5
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1376.5
$ws.Range("N110").Value = -5090

# Row 122 (Item ID 36168)
$ws.Range("H122").Value = 1604649
$ws.Range("I122").Value = 3664705
$ws.Range("J122").Value = 2383.3333
$ws.Range("K122").Value = 10994115
$ws.Range("L122").Value = 7149.999899999999
$ws.Range("M122").Value = -10991665
$ws.Range("N122").Value = -12049.9999

# Row 132 (Item ID 43997)
$ws.Range("H132").Value = 2099.3777
$ws.Range("I132").Value = 1855.3334
$ws.Range("J132").Value = 3685.6667
$ws.Range("K132").Value = 5566.0002
$ws.Range("L132").Value = 11057.0001
$ws.Range("M132").Value = -3036.0002
$ws.Range("N132").Value = -16117.0001

# Row 133 (Item ID 41857)
$ws.Range("H133").Value = 35964.668
$ws.Range("J133").Value = 35964.668
$ws.Range("L133").Value = 35964.668
$ws.Range("N133").Value = -41024.668

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (Item ID 19939)
$ws.Range("H94").Value = 2285.389
$ws.Range("I94").Value = 1484.75
$ws.Range("J94").Value = 2925.9
$ws.Range("K94").Value = 1484.75
$ws.Range("L94").Value = 2925.9
$ws.Range("M94").Value = -1033.75
$ws.Range("N94").Value = -3827.9

# Row 138 (Item ID 42308)
$ws.Range("H138").Value = 50696
$ws.Range("J138").Value = 50696
$ws.Range("L138").Value = 50696
$ws.Range("N138").Value = -60976

$ws = $wb.Worksheets.Item("CRP")
# Row 17 (Item ID 1823)
$ws.Range("H17").Value = 166694830
$ws.Range("J17").Value = 166694830
$ws.Range("L17").Value = 166694830
$ws.Range("N17").Value = -166695178

# Row 25 (Item ID 1895)
$ws.Range("H25").Value = 1000000000
$ws.Range("J25").Value = 1000000000
$ws.Range("L25").Value = 1000000000
$ws.Range("N25").Value = -1000000348

# Row 31 (Item ID 44023)
$ws.Range("H31").Value = 3273.86
$ws.Range("I31").Value = 2624.5
$ws.Range("J31").Value = 3784.0715
$ws.Range("K31").Value = 2624.5
$ws.Range("L31").Value = 3784.0715
$ws.Range("M31").Value = -2329.5
$ws.Range("N31").Value = -4374.0715

# Row 34 (Item ID 44023)
$ws.Range("H34").Value = 3273.86
$ws.Range("I34").Value = 2624.5
$ws.Range("J34").Value = 3784.0715
$ws.Range("K34").Value = 2624.5
$ws.Range("L34").Value = 3784.0715
$ws.Range("M34").Value = -2422.5
$ws.Range("N34").Value = -4188.0715

# Row 59 (Item ID 1942)
$ws.Range("H59").Value = 25895.7
$ws.Range("J59").Value = 25895.7
$ws.Range("L59").Value = 25895.7
$ws.Range("N59").Value = -28185.7

# Row 105 (Item ID 19928)
$ws.Range("H105").Value = 2255.5
$ws.Range("I105").Value = 1555
$ws.Range("J105").Value = 2956
$ws.Range("K105").Value = 1555
$ws.Range("L105").Value = 2956
$ws.Range("M105").Value = 192
$ws.Range("N105").Value = -6450

# Row 107 (Item ID 27689)
$ws.Range("H107").Value = 704.2
$ws.Range("I107").Value = 525.8570999999999
$ws.Range("J107").Value = 800.2308
$ws.Range("K107").Value = 525.8570999999999
$ws.Range("L107").Value = 800.2308
$ws.Range("M107").Value = 1394.1429
$ws.Range("N107").Value = -4640.2308

# Row 132 (Item ID 44019)
$ws.Range("H132").Value = 2210.7144
$ws.Range("I132").Value = 2210.7144
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6632.1432
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4102.1432
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 9 (Item ID 4681)
$ws.Range("H9").Value = 500000480
$ws.Range("I9").Value = 333333980
$ws.Range("J9").Value = 1000000000
$ws.Range("K9").Value = 1000001940
$ws.Range("L9").Value = 3000000000
$ws.Range("M9").Value = -1000001716
$ws.Range("N9").Value = -3000000448

# Row 74 (Item ID 12859)
$ws.Range("H74").Value = 17000
$ws.Range("J74").Value = 17000
$ws.Range("L74").Value = 51000
$ws.Range("N74").Value = -53122

# Row 77 (Item ID 12859)
$ws.Range("H77").Value = 17000
$ws.Range("J77").Value = 17000
$ws.Range("L77").Value = 153000
$ws.Range("N77").Value = -163608

# Row 141 (Item ID 44076)
$ws.Range("H141").Value = 13157.571
$ws.Range("I141").Value = 9908.385
$ws.Range("J141").Value = 18437.5
$ws.Range("K141").Value = 29725.155
$ws.Range("L141").Value = 55312.5
$ws.Range("M141").Value = -24545.155
$ws.Range("N141").Value = -65672.5

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (Item ID 36169)
$ws.Range("H102").Value = 831.6667
$ws.Range("I102").Value = 585.625
$ws.Range("K102").Value = 585.625
$ws.Range("M102").Value = 1036.375

# Row 122 (Item ID 36182)
$ws.Range("H122").Value = 3413962.5
$ws.Range("I122").Value = 4988253.5
$ws.Range("J122").Value = 2998.3333
$ws.Range("K122").Value = 14964760.5
$ws.Range("L122").Value = 8994.999899999999
$ws.Range("M122").Value = -14962310.5
$ws.Range("N122").Value = -13894.9999

# Row 126 (Item ID 36184)
$ws.Range("H126").Value = 12750.667
$ws.Range("I126").Value = 12750.667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 38252.001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -35782.001
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 136 (Item ID 44060)
$ws.Range("H136").Value = 4138.0312
$ws.Range("I136").Value = 2525.4644
$ws.Range("J136").Value = 15426
$ws.Range("K136").Value = 7576.3932
$ws.Range("L136").Value = 46278
$ws.Range("M136").Value = -5026.3932
$ws.Range("N136").Value = -51378

$ws = $wb.Worksheets.Item("WVR")
# Row 126 (Item ID 36210)
$ws.Range("H126").Value = 1468.8334
$ws.Range("I126").Value = 1202
$ws.Range("K126").Value = 3606
$ws.Range("M126").Value = -1136
